$wb = $excel.ActiveWorkbook

# --- Lvl0 (sheet1): add the "seasonEmptyField" object field definition ---
$ws0 = $wb.Worksheets.Item("Lvl0")
$ws0.Range("D19").Value = "field7"
$ws0.Range("E19").Value = "ref"
$ws0.Range("F19").Value = "seasonEmptyField"

# --- Lvl1 (sheet2): add the empty object itself (no fields) ---
$ws1 = $wb.Worksheets.Item("Lvl1")
$ws1.Range("C23").Value = "seasonEmptyField"

# --- Update selections to match where the author left the cursor ---
$ws0.Activate()
$ws0.Range("D20").Select()

$ws1.Activate()
$ws1.Range("C25").Select()
